# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 00:15"

# --- Simple statistic refreshes (country stays on the same row) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5835613
$ws.Range("C4").Value = 38886
$ws.Range("D4").Value = 3143507
$ws.Range("E4").Value = 2512029
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 877
$ws.Range("H4").Value = 180077

# Row 5: Brasil
$ws.Range("B5").Value = 3582362
$ws.Range("C5").Value = 45874
$ws.Range("D5").Value = 2670755
$ws.Range("E5").Value = 797357
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 796
$ws.Range("H5").Value = 114250

# Row 34: Egipto
$ws.Range("B34").Value = 97237
$ws.Range("C34").Value = 89
$ws.Range("D34").Value = 65118
$ws.Range("E34").Value = 26876
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 5243

# Row 48: Japon
$ws.Range("B48").Value = 60733
$ws.Range("C48").Value = 1012
$ws.Range("D48").Value = 47622
$ws.Range("E48").Value = 11942
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 1169

# Row 54: Barein
$ws.Range("B54").Value = 49038
$ws.Range("C54").Value = 377
$ws.Range("D54").Value = 45589
$ws.Range("E54").Value = 3266

# Row 148: Republica de Chipre
$ws.Range("B148").Value = 1417
$ws.Range("C148").Value = 11
$ws.Range("E148").Value = 519

# Row 168: Belice
$ws.Range("B168").Value = 668
$ws.Range("C168").Value = 20
$ws.Range("D168").Value = 44
$ws.Range("E168").Value = 618
$ws.Range("G168").Value = 1
$ws.Range("H168").Value = 6

# --- Reorder: Guyana moves up, immediately after Republica del Chad,
#     pushing Trinidad y Tobago and Santo Tome y Principe down one row ---

# Row 162 becomes Guyana (previously Trinidad y Tobago)
$ws.Range("A162").Value = "Guyana"
$ws.Range("B162").Value = 925
$ws.Range("C162").Value = 44
$ws.Range("D162").Value = 433
$ws.Range("E162").Value = 461
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 1
$ws.Range("H162").Value = 31

# Row 163 becomes Trinidad y Tobago (previously Santo Tome y Principe)
$ws.Range("A163").Value = "Trinidad yTobago"
$ws.Range("B163").Value = 899
$ws.Range("C163").Value = 35
$ws.Range("D163").Value = 165
$ws.Range("E163").Value = 721
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 13

# Row 164 becomes Santo Tome y Principe (previously Guyana)
$ws.Range("A164").Value = "Santo Tome y Principe"
$ws.Range("B164").Value = 892
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 831
$ws.Range("E164").Value = 46
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 15
